$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 771.1429000000001
$ws.Range("J19").Value = 951.55554
$ws.Range("L19").Value = 951.55554
$ws.Range("N19").Value = -1301.55554
$ws.Range("H40").Value = 10461.75
$ws.Range("I40").Value = 3845
$ws.Range("J40").Value = 12667.333
$ws.Range("K40").Value = 3845
$ws.Range("L40").Value = 12667.333
$ws.Range("M40").Value = -3670
$ws.Range("N40").Value = -13017.333
$ws.Range("H88").Value = 5798.2
$ws.Range("I88").Value = 2497
$ws.Range("J88").Value = 7999
$ws.Range("K88").Value = 2497
$ws.Range("L88").Value = 7999
$ws.Range("M88").Value = -2091
$ws.Range("N88").Value = -8811
$ws.Range("H91").Value = 5798.2
$ws.Range("I91").Value = 2497
$ws.Range("J91").Value = 7999
$ws.Range("K91").Value = 2497
$ws.Range("L91").Value = 7999
$ws.Range("M91").Value = -1093
$ws.Range("N91").Value = -10807
$ws.Range("H96").Value = 544.8889
$ws.Range("I96").Value = 267.5
$ws.Range("K96").Value = 802.5
$ws.Range("M96").Value = 570.5
$ws.Range("H138").Value = 2884.2708
$ws.Range("J138").Value = 3639.8
$ws.Range("L138").Value = 10919.4
$ws.Range("N138").Value = -21199.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4003.5386
$ws.Range("I32").Value = 3385.2942
$ws.Range("K32").Value = 3385.2942
$ws.Range("M32").Value = -3098.2942
$ws.Range("H74").Value = 27784260
$ws.Range("I74").Value = 37043336
$ws.Range("J74").Value = 7033.3335
$ws.Range("K74").Value = 37043336
$ws.Range("L74").Value = 7033.3335
$ws.Range("M74").Value = -37042462
$ws.Range("N74").Value = -8781.333500000001
$ws.Range("H77").Value = 27784260
$ws.Range("I77").Value = 37043336
$ws.Range("J77").Value = 7033.3335
$ws.Range("K77").Value = 185216680
$ws.Range("L77").Value = 35166.6675
$ws.Range("M77").Value = -185212312
$ws.Range("N77").Value = -43902.6675
$ws.Range("H110").Value = 6976.364
$ws.Range("I110").Value = 6214.1113
$ws.Range("K110").Value = 6214.1113
$ws.Range("M110").Value = -4169.1113
$ws.Range("H122").Value = 2929.8
$ws.Range("I122").Value = 2090.238
$ws.Range("K122").Value = 6270.714
$ws.Range("M122").Value = -3820.714
$ws.Range("H132").Value = 8059.727
$ws.Range("I132").Value = 4088.7144
$ws.Range("K132").Value = 12266.1432
$ws.Range("M132").Value = -9736.143199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 503.5357
$ws.Range("I80").Value = 441.33334
$ws.Range("J80").Value = 533
$ws.Range("K80").Value = 441.33334
$ws.Range("L80").Value = 533
$ws.Range("M80").Value = 556.66666
$ws.Range("N80").Value = -2529
$ws.Range("H83").Value = 503.5357
$ws.Range("I83").Value = 441.33334
$ws.Range("J83").Value = 533
$ws.Range("K83").Value = 2206.6667
$ws.Range("L83").Value = 2665
$ws.Range("M83").Value = 2785.3333
$ws.Range("N83").Value = -12649
$ws.Range("H86").Value = 3011.8096
$ws.Range("I86").Value = 2071.7222
$ws.Range("J86").Value = 8652.333000000001
$ws.Range("K86").Value = 2071.7222
$ws.Range("L86").Value = 8652.333000000001
$ws.Range("M86").Value = -948.7222000000002
$ws.Range("N86").Value = -10898.333
$ws.Range("H89").Value = 3011.8096
$ws.Range("I89").Value = 2071.7222
$ws.Range("J89").Value = 8652.333000000001
$ws.Range("K89").Value = 10358.611
$ws.Range("L89").Value = 43261.665
$ws.Range("M89").Value = -4742.611000000001
$ws.Range("N89").Value = -54493.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1587.8182
$ws.Range("I107").Value = 2253.6667
$ws.Range("K107").Value = 2253.6667
$ws.Range("M107").Value = -333.6667000000002
$ws.Range("H132").Value = 5397.5454
$ws.Range("I132").Value = 4624.4287
$ws.Range("J132").Value = 6750.5
$ws.Range("K132").Value = 13873.2861
$ws.Range("L132").Value = 20251.5
$ws.Range("M132").Value = -11343.2861
$ws.Range("N132").Value = -25311.5
$ws.Range("H134").Value = 7965.1113
$ws.Range("I134").Value = 4060.5881
$ws.Range("J134").Value = 14602.8
$ws.Range("K134").Value = 12181.7643
$ws.Range("L134").Value = 43808.39999999999
$ws.Range("M134").Value = -9646.764299999999
$ws.Range("N134").Value = -48878.39999999999
$ws.Range("H141").Value = 237191.6
$ws.Range("J141").Value = 237191.6
$ws.Range("L141").Value = 237191.6
$ws.Range("N141").Value = -247551.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 110.25
$ws.Range("I7").Value = 94.333336
$ws.Range("J7").Value = 126.166664
$ws.Range("K7").Value = 283.000008
$ws.Range("L7").Value = 378.499992
$ws.Range("M7").Value = -171.000008
$ws.Range("N7").Value = -602.499992
$ws.Range("H113").Value = 1213.5454
$ws.Range("I113").Value = 1209.2
$ws.Range("J113").Value = 1217.1666
$ws.Range("K113").Value = 3627.6
$ws.Range("L113").Value = 3651.4998
$ws.Range("M113").Value = -1457.6
$ws.Range("N113").Value = -7991.4998
$ws.Range("H122").Value = 1514.6818
$ws.Range("J122").Value = 1707.7368
$ws.Range("L122").Value = 15369.6312
$ws.Range("N122").Value = -20269.6312
$ws.Range("H129").Value = 5954466
$ws.Range("I129").Value = 832.36365
$ws.Range("K129").Value = 2497.09095
$ws.Range("M129").Value = 2502.90905
$ws.Range("H138").Value = 5547
$ws.Range("J138").Value = 13265.333
$ws.Range("L138").Value = 39795.999
$ws.Range("N138").Value = -50075.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4101.75
$ws.Range("I22").Value = 4908
$ws.Range("J22").Value = 3833
$ws.Range("K22").Value = 4908
$ws.Range("L22").Value = 3833
$ws.Range("M22").Value = -4379
$ws.Range("N22").Value = -4891
$ws.Range("H70").Value = 9999.666999999999
$ws.Range("I70").Value = 9999
$ws.Range("K70").Value = 9999
$ws.Range("M70").Value = -9729
$ws.Range("H73").Value = 9999.666999999999
$ws.Range("I73").Value = 9999
$ws.Range("K73").Value = 9999
$ws.Range("M73").Value = -9063
$ws.Range("H80").Value = 7572.2856
$ws.Range("I80").Value = 5666.6665
$ws.Range("K80").Value = 5666.6665
$ws.Range("M80").Value = -4668.6665
$ws.Range("H83").Value = 7572.2856
$ws.Range("I83").Value = 5666.6665
$ws.Range("K83").Value = 28333.3325
$ws.Range("M83").Value = -23341.3325
$ws.Range("H97").Value = 418
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 224
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 224
$ws.Range("M97").Value = -504
$ws.Range("N97").Value = -1216
$ws.Range("H122").Value = 4470.4443
$ws.Range("I122").Value = 3230.6667
$ws.Range("K122").Value = 9692.000100000001
$ws.Range("M122").Value = -7242.000100000001
$ws.Range("H141").Value = 57996.332
$ws.Range("J141").Value = 57996.332
$ws.Range("L141").Value = 57996.332
$ws.Range("N141").Value = -68356.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15416.5
$ws.Range("J22").Value = 15416.5
$ws.Range("L22").Value = 15416.5
$ws.Range("N22").Value = -16006.5
$ws.Range("H27").Value = 15416.5
$ws.Range("J27").Value = 15416.5
$ws.Range("L27").Value = 15416.5
$ws.Range("N27").Value = -15630.5
$ws.Range("H40").Value = 9625.066000000001
$ws.Range("I40").Value = 7047.1
$ws.Range("K40").Value = 7047.1
$ws.Range("M40").Value = -6911.1
$ws.Range("H93").Value = 64999
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 111993.445
$ws.Range("J141").Value = 111993.445
$ws.Range("L141").Value = 111993.445
$ws.Range("N141").Value = -122353.445
